# Append new 2h OHLCV rows (1169-1232) to LINKUSDT_2h sheet
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$startRow = 1169

$data = @(
    @(45534.5, 11.03, 11.17, 10.95, 11.06, 142813.45),
    @(45534.58333333334, 11.07, 11.12, 10.51, 10.54, 407614.3),
    @(45534.66666666666, 10.54, 10.86, 10.5, 10.82, 218387.27),
    @(45534.75, 10.82, 11.11, 10.81, 11, 149530.47),
    @(45534.83333333334, 10.99, 11.08, 10.96, 11.02, 81317.82000000001),
    @(45534.91666666666, 11.03, 11.14, 11.02, 11.11, 81950.06),
    @(45535, 11.12, 11.18, 11.09, 11.11, 85481.49000000001),
    @(45535.08333333334, 11.1, 11.14, 11.07, 11.11, 40008.88),
    @(45535.16666666666, 11.12, 11.18, 11.1, 11.13, 56387.84),
    @(45535.25, 11.14, 11.14, 11.09, 11.12, 30943.54),
    @(45535.33333333334, 11.13, 11.16, 11.03, 11.06, 77054.12),
    @(45535.41666666666, 11.07, 11.2, 11.05, 11.15, 81948.03),
    @(45535.5, 11.14, 11.21, 11.08, 11.08, 68943.63),
    @(45535.58333333334, 11.09, 11.12, 11.03, 11.09, 62844.85),
    @(45535.66666666666, 11.09, 11.11, 10.93, 10.99, 106138.48),
    @(45535.75, 10.99, 11.06, 10.98, 11.06, 14887.74),
    @(45535.83333333334, 11.02, 11.06, 10.99, 11.06, 28981.75),
    @(45535.91666666666, 11.07, 11.08, 11, 11.02, 40196.9),
    @(45536, 11.02, 11.02, 10.9, 10.91, 64567.54),
    @(45536.08333333334, 10.9, 10.98, 10.83, 10.94, 100645.15),
    @(45536.16666666666, 10.94, 10.97, 10.78, 10.8, 105162.46),
    @(45536.25, 10.81, 10.9, 10.79, 10.88, 81507.22),
    @(45536.33333333334, 10.88, 10.9, 10.82, 10.85, 41717.51),
    @(45536.41666666666, 10.85, 10.85, 10.81, 10.81, 8198.959999999999),
    @(45536.5, 10.8, 10.82, 10.61, 10.72, 138340.44),
    @(45536.58333333334, 10.71, 10.75, 10.52, 10.7, 210150.84),
    @(45536.66666666666, 10.7, 10.74, 10.52, 10.6, 142488.58),
    @(45536.75, 10.61, 10.77, 10.58, 10.71, 114132.14),
    @(45536.83333333334, 10.71, 10.77, 10.55, 10.62, 106420.51),
    @(45536.91666666666, 10.62, 10.63, 10.23, 10.37, 280412.5),
    @(45537, 10.37, 10.49, 10.32, 10.44, 136245.87),
    @(45537.08333333334, 10.43, 10.51, 10.4, 10.48, 54566.67),
    @(45537.16666666666, 10.48, 10.52, 10.34, 10.48, 95180.46000000001),
    @(45537.25, 10.49, 10.52, 10.33, 10.38, 76160.02),
    @(45537.33333333334, 10.38, 10.7, 10.37, 10.68, 148621.65),
    @(45537.41666666666, 10.68, 10.79, 10.67, 10.75, 146799.14),
    @(45537.5, 10.75, 10.76, 10.65, 10.68, 113864.09),
    @(45537.58333333334, 10.67, 10.79, 10.59, 10.66, 194184.65),
    @(45537.66666666666, 10.66, 10.69, 10.55, 10.67, 128867.89),
    @(45537.75, 10.68, 10.71, 10.62, 10.66, 61154.97),
    @(45537.83333333334, 10.67, 10.84, 10.64, 10.77, 113839.63),
    @(45537.91666666666, 10.78, 10.8, 10.73, 10.77, 45473.72),
    @(45538, 10.77, 10.82, 10.73, 10.8, 70268.86),
    @(45538.08333333334, 10.79, 10.88, 10.69, 10.7, 88114.36),
    @(45538.16666666666, 10.7, 10.73, 10.65, 10.71, 77782.00999999999),
    @(45538.25, 10.7, 10.71, 10.61, 10.69, 96711.42999999999),
    @(45538.33333333334, 10.7, 10.71, 10.54, 10.56, 86943.5),
    @(45538.41666666666, 10.56, 10.67, 10.55, 10.67, 99072.08),
    @(45538.5, 10.67, 10.74, 10.41, 10.47, 221183.21),
    @(45538.58333333334, 10.46, 10.46, 10.27, 10.3, 218753.58),
    @(45538.66666666666, 10.29, 10.35, 10.24, 10.25, 142557.4),
    @(45538.75, 10.26, 10.51, 10.21, 10.5, 246966.18),
    @(45538.83333333334, 10.49, 10.52, 10.36, 10.4, 126456.66),
    @(45538.91666666666, 10.4, 10.42, 10.29, 10.31, 121924.97),
    @(45539, 10.32, 10.39, 9.82, 10.19, 542025.48),
    @(45539.08333333334, 10.19, 10.31, 10.15, 10.25, 155309.17),
    @(45539.16666666666, 10.25, 10.31, 10.17, 10.19, 104417.88),
    @(45539.25, 10.2, 10.41, 10.2, 10.26, 208190.12),
    @(45539.33333333334, 10.27, 10.3, 10.22, 10.22, 97844.3),
    @(45539.41666666666, 10.22, 10.24, 10.09, 10.12, 143167.84),
    @(45539.5, 10.12, 10.18, 10.06, 10.11, 171334.26),
    @(45539.58333333334, 10.11, 10.41, 10.06, 10.4, 390478.22),
    @(45539.66666666666, 10.39, 10.5, 10.33, 10.34, 190172.8),
    @(45539.75, 10.34, 10.39, 10.27, 10.29, 126129.72)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $startRow + $i
    $vals = $data[$i]
    $ws.Cells.Item($r, 1).Value = $vals[0]
    $ws.Cells.Item($r, 2).Value = $vals[1]
    $ws.Cells.Item($r, 3).Value = $vals[2]
    $ws.Cells.Item($r, 4).Value = $vals[3]
    $ws.Cells.Item($r, 5).Value = $vals[4]
    $ws.Cells.Item($r, 6).Value = $vals[5]
}

# Re-apply the timestamp column date style (same as existing rows) to the new rows
$endRow = $startRow + $data.Count - 1
$srcDateCell = $ws.Range("A1168")
$dstDateRange = $ws.Range("A" + $startRow + ":A" + $endRow)
$srcDateCell.Copy()
$dstDateRange.PasteSpecial(-4122)
$excel.CutCopyMode = 0

